# "working rate of fire" -- Events.xlsx edit
#
# 1. Insert a new table row right after the "e052 Pivot Tank" entry for a
#    new "e052a" / "Pivot Turret" entry (lets the player rotate the turret
#    independently of the hull).
# 2. Give the "e053c Main Gun Firing - Rate of Fire" entry's body text two
#    extra trailing line breaks so it lays out consistently with its
#    neighbours.
# 3. Leave the selection on the row that was being edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the short event code, column B the rich "event body" text.
# Style 1 (vertical-top) is used for column A, style 2 (left/top + wrap,
# text format) is used for column B -- mirror both on the new row.
$xlTop = -4160      # xlTop
$xlLeft = -4131     # xlLeft

# --- 1. Insert the new "e052a" / Pivot Turret row right after row 62 (e052) ---
$ws.Rows(63).Insert()

$ws.Range("A63").Value = "e052a"
$ws.Range("A63").VerticalAlignment = $xlTop

$pivotTurretText = @"
<Bold>e052 Pivot Turret</Bold> 
<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r8.24' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
Select the plus or minus buttons to rotate.  
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Button Content='   -   ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<InlineUIContainer><Image Name='ShermanPivot'  Height='150' Width='150'></Image></InlineUIContainer> 
<InlineUIContainer><Button Content='   +   ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
When you are satisfied with the current orientation, click Sherman image between buttons to continue.
"@
$ws.Range("B63").Value = $pivotTurretText
$ws.Range("B63").VerticalAlignment = $xlTop
$ws.Range("B63").HorizontalAlignment = $xlLeft
$ws.Range("B63").WrapText = $true
$ws.Range("B63").NumberFormat = "@"

$ws.Rows(63).RowHeight = 158.4

# --- 2. Append trailing line breaks to the "Rate of Fire" entry (now row 67) ---
$rateOfFireText = @"
<Bold>e053c Main Gun Firing - Rate of Fire</Bold> 
<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r9.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
Rolled low enough on the 
<InlineUIContainer><Button Content='Rate of Fire' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table to repeat fire on same target.
<LineBreak/><LineBreak/>
 Choose 
 <InlineUIContainer><Button Content='Fire' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   or  
 <InlineUIContainer><Button Content='Skip' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   to continue.
<LineBreak/><LineBreak/>
"@
$ws.Range("B67").Value = $rateOfFireText
$ws.Rows(67).RowHeight = 158.4

# --- 3. Update the active selection to match where editing left off ---
$ws.Range("B67").Select()
